$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for rows with refreshed data ---
$ws.Range("D2").Value = "29.474.24"
$ws.Range("E2").Value = "  -3.03%  "
$ws.Range("D3").Value = "1.993.87"
$ws.Range("E3").Value = "  -4.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.016"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.014"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5025"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4234"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.77"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08939"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.111"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.03%  "
$ws.Range("E12").Value = "  -6.32%  "
$ws.Range("D13").Value = "2.002.77"
$ws.Range("E13").Value = "  -7.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.963"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.455"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.61%  "
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001112"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06775"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.014"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.936"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.35%  "
$ws.Range("D23").Value = "29.522.69"
$ws.Range("E23").Value = "  -3.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.314"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.292"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.303"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.059"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09942"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.543"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.829"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.793"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02457"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.236"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06388"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.16%  "
$ws.Range("E39").Value = "  -3.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6529"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2041"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.18%  "
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6326"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.309"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.505"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000340"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.04%  "
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.131"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.74%  "

# --- Rows 45/46: EnergySwap and NEARProtocol swapped positions with refreshed data ---
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.07%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.211"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.66%  "
